# "Generate Report for Handback"
#
# A new handback run finished for the b47a380c-... file: its handoff/handback
# timestamps got refreshed (everything else - zh-cn row 2, de-de row 2, the
# 271e6631-... rows - stayed "in sync" and unchanged).
#
#   Overview!G3  (Latest HO Xliff Generate Date)   05:01:40 -> 05:02:51
#   zh-cn!H3     (Correspond Handoff Datetime)     05:01:34 -> 05:02:45
#   zh-cn!K3     (Correspond Handback DateTime)    05:02:21 -> 05:03:10
#   de-de!H3     (Correspond Handoff Datetime)     05:01:40 -> 05:02:51
#   de-de!K3     (Correspond Handback DateTime)    05:02:29 -> 05:03:19

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-07 05:02:51"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-07 05:02:45"
$wsZhCn.Range("K3").Value = "2016-09-07 05:03:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-07 05:02:51"
$wsDeDe.Range("K3").Value = "2016-09-07 05:03:19"
